$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41 / 42 content swap (ApeXProtocol <-> Fetch.AI) ---
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.63'
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.39'
$ws.Range("E42").Value = '  +0.06%  '

# --- Remaining per-cell price / volume updates ---
$ws.Range("D2").Value = '67.570.17'
$ws.Range("E2").Value = '  +2.94%  '
$ws.Range("D3").Value = '3.271.15'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.33'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.53'
$ws.Range("E6").Value = '  -2.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.583'
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("D9").Value = '3.264.68'
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.176'
$ws.Range("E10").Value = '  +0.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.572'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.38'
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000269'
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '682.87'
$ws.Range("E14").Value = '  +11.56%  '
$ws.Range("D15").Value = '3.790.13'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.34'
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").Value = '67.586.99'
$ws.Range("E17").Value = '  +2.97%  '
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").Value = '3.264.31'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.32'
$ws.Range("E20").Value = '  -2.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.71'
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.889'
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.97'
$ws.Range("E23").Value = '  -5.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.14'
$ws.Range("E24").Value = '  +3.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.02'
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.34'
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.49'
$ws.Range("E29").Value = '  +5.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.43'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.64'
$ws.Range("E31").Value = '  +3.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '580.94'
$ws.Range("E32").Value = '  +6.18%  '
$ws.Range("D33").Value = '3.860.97'
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.80'
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.36'
$ws.Range("E37").Value = '  -9.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.33'
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.130'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("E40").Value = '  +2.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '32.05'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("D44").Value = '0.0₃0670'
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.329'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.39'
$ws.Range("E49").Value = '  +9.65%  '
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.03'
$ws.Range("E51").Value = '  +0.69%  '
